$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of test data (rows 121 and 122), in the order that
# reproduces the original author's shared-string insertion order.
$ws.Range("C122").Value = "UpsidePotentialRatio_test2"
$ws.Range("B122").Value = "Test upside potential ratio for full sets"
$ws.Range("A122").Value = "UpsidePotentialRatio2"

$ws.Range("A121").Value = "UpsidePotentialRatio1"
$ws.Range("B121").Value = "Test upside potential ratio for subsets"
$ws.Range("C121").Value = "UpsidePotentialRatio_test1"

$ws.Range("C121").Select()
